# Update the "取得日時" (retrieved datetime) column (A) for rows 2-9 on the
# "ランサーズ" sheet, replacing the old timestamp with the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-25 06:24:40"
$newTimestamp = "2025-10-25 06:36:18"

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
